$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the date-looking string into A86 as literal text (matching the
# existing rows, which store dates as plain text rather than date serials),
# then restore the default (unstyled) cell formatting so the new row looks
# just like its neighbors.
$ws.Range("A86").NumberFormat = "@"
$ws.Range("A86").Value = "11/11/2025"
$ws.Range("A86").Style = "Normal"

$ws.Range("B86").Value = 10269.31
